$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = '38.32***'
$ws.Range("E3").Value = '[33.52, 43.92]'
$ws.Range("F3").Value = '110.62***'
$ws.Range("G3").Value = '[99.72, 122.78]'
$ws.Range("H3").Value = '3.66***'
$ws.Range("I3").Value = '[ 3.45, 3.86]'
$ws.Range("B5").Value = '1.53***'
$ws.Range("C5").Value = '[ 1.35,  1.77]'
$ws.Range("D5").Value = '1.03'
$ws.Range("E5").Value = '[ 0.98,  1.08]'
$ws.Range("K5").Value = '[0.70,   1.02]'
$ws.Range("C6").Value = '[ 1.22,  1.54]'
$ws.Range("D6").Value = '1.03'
$ws.Range("C7").Value = '[ 0.73,  1.31]'
$ws.Range("E7").Value = '[ 0.82,  0.99]'
$ws.Range("G7").Value = '[ 0.89,   1.02]'
$ws.Range("I7").Value = '[-0.13, 0.07]'
$ws.Range("K7").Value = '[1.19,   4.16]'
$ws.Range("C8").Value = '[ 1.13,  2.46]'
$ws.Range("I8").Value = '[-0.10, 0.10]'
$ws.Range("K8").Value = '[0.60,   3.88]'
$ws.Range("C11").Value = '[ 0.64,  1.14]'
$ws.Range("I11").Value = '[ 0.15, 0.36]'
$ws.Range("J11").Value = '1.87'
$ws.Range("K11").Value = '[0.89,   3.84]'
$ws.Range("B12").Value = '8.93***'
$ws.Range("C12").Value = '[ 7.39, 10.84]'
$ws.Range("E12").Value = '[ 1.21,  1.46]'
$ws.Range("J12").Value = '1.17'
$ws.Range("K12").Value = '[0.66,   2.07]'
$ws.Range("C13").Value = '[ 1.02,  1.48]'
$ws.Range("D13").Value = '1.07'
$ws.Range("E13").Value = '[ 0.99,  1.16]'
$ws.Range("J13").Value = '0.87'
$ws.Range("K13").Value = '[0.51,   1.49]'
$ws.Range("B16").Value = '1.59'
$ws.Range("C16").Value = '[ 0.83,  2.99]'
$ws.Range("D16").Value = '1.07'
$ws.Range("E16").Value = '[ 0.81,  1.39]'
$ws.Range("G16").Value = '[ 0.86,   1.40]'
$ws.Range("I16").Value = '[-0.09, 0.93]'
$ws.Range("J16").Value = '1.89'
$ws.Range("K16").Value = '[0.62,   6.10]'
$ws.Range("B17").Value = '1.64'
$ws.Range("C17").Value = '[ 0.85,  3.10]'
$ws.Range("D17").Value = '1.04'
$ws.Range("E17").Value = '[ 0.79,  1.35]'
$ws.Range("F17").Value = '1.02'
$ws.Range("G17").Value = '[ 0.80,   1.32]'
$ws.Range("I17").Value = '[-0.17, 0.85]'
$ws.Range("K17").Value = '[0.32,   3.64]'
$ws.Range("B18").Value = '0.33**'
$ws.Range("C18").Value = '[ 0.15,  0.71]'
$ws.Range("D18").Value = '1.19'
$ws.Range("E18").Value = '[ 0.86,  1.65]'
$ws.Range("F18").Value = '0.91'
$ws.Range("G18").Value = '[ 0.69,   1.20]'
$ws.Range("H18").Value = '-0.36'
$ws.Range("I18").Value = '[-0.90, 0.19]'
$ws.Range("J18").Value = '18.19***'
$ws.Range("K18").Value = '[3.23, 126.21]'
$ws.Range("B19").Value = '0.55'
$ws.Range("C19").Value = '[ 0.26,  1.20]'
$ws.Range("D19").Value = '0.99'
$ws.Range("E19").Value = '[ 0.71,  1.37]'
$ws.Range("F19").Value = '1.04'
$ws.Range("G19").Value = '[ 0.79,   1.36]'
$ws.Range("I19").Value = '[-0.83, 0.27]'
$ws.Range("J19").Value = '1.19'
$ws.Range("K19").Value = '[0.16,   7.82]'
$ws.Range("B24").Value = '0.72'
$ws.Range("C24").Value = '[0.55, 0.97]'
$ws.Range("E24").Value = '[0.23, 0.39]'
$ws.Range("I24").Value = '[0.48, 0.78]'
$ws.Range("K24").Value = '[0.81, 1.80]'
$ws.Range("C25").Value = '[0.05, 0.40]'
$ws.Range("D25").Value = '0.11'
$ws.Range("E25").Value = '[0.07, 0.17]'
$ws.Range("K25").Value = '[0.01, 0.50]'
$ws.Range("B26").Value = '0.14'
$ws.Range("C26").Value = '[0.01, 0.33]'
$ws.Range("E26").Value = '[0.04, 0.13]'
$ws.Range("J26").Value = '0.44'
$ws.Range("B27").Value = '0.18'
$ws.Range("C27").Value = '[0.01, 0.69]'
$ws.Range("K27").Value = '[0.05, 2.09]'
$ws.Range("B28").Value = '0.24'
$ws.Range("C28").Value = '[0.01, 0.91]'
$ws.Range("D28").Value = '0.06'
$ws.Range("E28").Value = '[0.00, 0.19]'
$ws.Range("J28").Value = '0.74'
$ws.Range("K28").Value = '[0.03, 2.66]'
$ws.Range("E32").Value = '[0.66, 0.71]'
